# Re-run of the backward-elimination export re-saved the 29 OLS summary
# sheets; each sheet's B2 cell holds the statsmodels ".summary()" text
# block whose "Date:" / "Time:" header line is stamped at render time.
# This refreshes those stamps to the new run (Sun, 05 Jan 2020, ~21:22)
# while leaving every numeric regression result untouched.

$wb = $excel.ActiveWorkbook

$oldDate = "Thu, 02 Jan 2020"
$newDate = "Sun, 05 Jan 2020"

# Per-sheet old -> new "Time:" stamp (sheet 1 = first tab ... sheet 29 = last tab).
$newTime1 = "21:22:09"
$newTime2 = "21:22:10"

$oldTimes = @(
    "20:48:31", "20:48:31", "20:48:31", "20:48:31", "20:48:31", "20:48:31",
    "20:48:31", "20:48:31", "20:48:31", "20:48:31", "20:48:31", "20:48:31",
    "20:48:32", "20:48:32", "20:48:32", "20:48:32", "20:48:32", "20:48:32",
    "20:48:32", "20:48:32", "20:48:32", "20:48:32", "20:48:32", "20:48:32",
    "20:48:32", "20:48:32", "20:48:32", "20:48:32", "20:48:32"
)

$newTimes = @(
    $newTime1, $newTime1, $newTime1, $newTime1, $newTime1, $newTime1,
    $newTime1, $newTime1, $newTime1, $newTime1, $newTime1, $newTime1,
    $newTime1, $newTime1, $newTime2, $newTime2, $newTime2, $newTime2,
    $newTime2, $newTime2, $newTime2, $newTime2, $newTime2, $newTime2,
    $newTime2, $newTime2, $newTime2, $newTime2, $newTime2
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $cell = $ws.Range("B2")
    $text = $cell.Text

    if ($text -and $text.Contains("Date:")) {
        $oldTime = $oldTimes[$i - 1]
        $newTime = $newTimes[$i - 1]

        # Remember the row's current height: re-assigning a long wrapped
        # string via .Value makes the host re-derive a wrap height, which
        # does not match the original file's stamped height. Restore it
        # below so only the cell text (the actual authored change) moves.
        $row = $ws.Rows.Item(2)
        $originalHeight = $row.RowHeight

        $updated = $text.Replace($oldDate, $newDate).Replace($oldTime, $newTime)
        $cell.Value = $updated

        $row.RowHeight = $originalHeight
    }
}
